$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 15.70730000000002
$ws.Range("E4").Value = 16.28759999999999
$ws.Range("C6").Value = -12.5006
$ws.Range("B7").Value = 5.638299999999999
$ws.Range("D7").Value = -7.666899999999994
$ws.Range("A9").Value = -21.95300000000002
$ws.Range("D10").Value = -8.361499999999999
$ws.Range("E11").Value = 15.5682
$ws.Range("B12").Value = 5.466799999999999
$ws.Range("A13").Value = -22.3342
$ws.Range("D13").Value = -8.461999999999996
$ws.Range("E13").Value = 16.09380000000001
$ws.Range("B14").Value = 5.563300000000001
$ws.Range("C15").Value = -14.38340000000001
$ws.Range("A16").Value = -21.52729999999998
$ws.Range("D16").Value = -9.123800000000006
$ws.Range("A18").Value = -22.32050000000001
$ws.Range("B19").Value = 8.832100000000004
$ws.Range("A20").Value = -20.66699999999998
$ws.Range("D20").Value = -6.883299999999997
$ws.Range("D24").Value = -7.6238
$ws.Range("E24").Value = 16.52750000000001
$ws.Range("E25").Value = 17.11830000000001
$ws.Range("A26").Value = -21.21349999999998
$ws.Range("B26").Value = 4.100600000000004
$ws.Range("A27").Value = -22.07109999999999
$ws.Range("B27").Value = 5.613600000000004
$ws.Range("E27").Value = 16.98449999999998
$ws.Range("C28").Value = -12.50049999999999
$ws.Range("A29").Value = -21.8244
$ws.Range("B29").Value = 5.631699999999999
$ws.Range("E29").Value = 17.347
$ws.Range("D32").Value = -9.016899999999996
$ws.Range("C33").Value = -11.9584
$ws.Range("A35").Value = -19.98980000000001
$ws.Range("C35").Value = -12.62430000000001
$ws.Range("E35").Value = 16.44489999999999
$ws.Range("A36").Value = -20.126
$ws.Range("B37").Value = 8.283600000000012
$ws.Range("B38").Value = 4.340000000000001
$ws.Range("C38").Value = -12.1576
$ws.Range("D39").Value = -7.230900000000005
$ws.Range("E40").Value = 17.08410000000002
$ws.Range("C43").Value = -14.3467
$ws.Range("C44").Value = -13.6376
$ws.Range("E44").Value = 16.41189999999998
$ws.Range("A45").Value = -22.15019999999999
$ws.Range("C45").Value = -13.3125
$ws.Range("B47").Value = 5.537999999999999
$ws.Range("C47").Value = -12.2717
$ws.Range("D47").Value = -7.250499999999995
$ws.Range("D48").Value = -7.069999999999998
$ws.Range("E49").Value = 16.47249999999999
$ws.Range("B51").Value = 6.105900000000005
$ws.Range("C51").Value = -11.9255
$ws.Range("B52").Value = 5.793000000000001
$ws.Range("D52").Value = -7.263899999999995
$ws.Range("C54").Value = -13.3739
$ws.Range("A55").Value = -22.399
$ws.Range("B55").Value = 4.574499999999996
$ws.Range("D56").Value = -7.832699999999999
$ws.Range("A57").Value = -22.0122
$ws.Range("C57").Value = -13.17779999999999
$ws.Range("E57").Value = 16.48639999999999
$ws.Range("C62").Value = -13.98310000000001
$ws.Range("C63").Value = -11.8066
$ws.Range("C67").Value = -10.4245
$ws.Range("A69").Value = -21.68090000000001
$ws.Range("B69").Value = 5.552499999999996
$ws.Range("B70").Value = 6.176000000000005
$ws.Range("C70").Value = -12.29819999999999
$ws.Range("A76").Value = -22.4064
$ws.Range("B76").Value = 5.6446
$ws.Range("A78").Value = -19.89799999999999
$ws.Range("E80").Value = 17.0776
$ws.Range("B81").Value = 5.756000000000001
$ws.Range("C81").Value = -12.5519
$ws.Range("A82").Value = -21.9825
$ws.Range("A83").Value = -21.74959999999999
$ws.Range("B83").Value = 6.246300000000009
$ws.Range("D84").Value = -8.565099999999994
$ws.Range("E85").Value = 16.2594
$ws.Range("C88").Value = -11.58719999999999
$ws.Range("E89").Value = 17.40180000000002
$ws.Range("A93").Value = -20.64029999999998
$ws.Range("B94").Value = 5.709799999999997
$ws.Range("C96").Value = -12.41780000000001
$ws.Range("A97").Value = -21.8776
$ws.Range("C99").Value = -12.4035
$ws.Range("B100").Value = 4.9581
$ws.Range("D100").Value = -8.305600000000004
$ws.Range("D101").Value = -8.006899999999996
$ws.Range("E101").Value = 16.4201
$ws.Range("B102").Value = 8.209500000000007
